# Generate Report for Handoff
# Update the "Latest Handoff Datetime" (column D, row 6 - the ca03d5d7 file)
# on both the zh-cn and de-de localization-status sheets with freshly
# generated handoff timestamps.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("D6").Value = "2016-03-08 08:24:54"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("D6").Value = "2016-03-08 08:24:58"
